# Editace rezervace, seřazení rezervací podle času
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new reservation row 31: task name (A31) and hours (B31)
$ws.Range("A31").Value = "4. iterace - další implementace (vybrat stůl)"
$ws.Range("B31").Value = 1

# Move selection to A32, matching the updated sheet view selection
$ws.Range("A32").Select()

$wb.Save()
